$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: merge the old "E" column (asdff/ddd) away, letting the old
# "F" column (Result/PAkSS) become the new, final "E" column, fixing the
# "PAkSS" typo to "PASS" in the process.
$ws2.Range("E1").Value = "Result"
$ws2.Range("E1").Style = "Normal"
$ws2.Range("E2").Value = "PASS"
$ws2.Range("F1:F2").Clear()

# --- Sheet2 B2: rename the old "testEasy" test case to "apptesting" and
# give it its own look: 9pt black Courier New, vertically centred.
$ws2.Range("B2").Value = "apptesting"
$ws2.Range("B2").Font.Name = "Courier New"
$ws2.Range("B2").Font.Size = 9
$ws2.Range("B2").Font.Color = 0
$ws2.Range("B2").VerticalAlignment = -4108

# --- Sheet2 column widths: column B now needs a best-fit custom width,
# matching the existing best-fit width already on column D.
$ws2.Columns.Item(2).ColumnWidth = 10.15

# --- Sheet2 view: selection moved, and it is no longer the active tab.
$ws2.Range("B4").Select() | Out-Null

# --- Sheet2 print setup: portrait orientation.
$ws2.PageSetup.Orientation = 1

# --- Sheet1 becomes the active / selected sheet in the workbook.
$ws1.Activate() | Out-Null
$ws1.Range("B3").Select() | Out-Null
